$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5203476666666668
$ws.Range("H2").Value = 1.561043
$ws.Range("I2").Value = 0.004105934376266647
$ws.Range("J2").Value = 0.004105934376266647
$ws.Range("M2").Value = 0.4435246666666666
$ws.Range("N2").Value = 1.330574
$ws.Range("O2").Value = 0.02616294045579069
$ws.Range("P2").Value = 0.02616294045579069
$ws.Range("Q2").Value = 0.2307870254091111
$ws.Range("R2").Value = 2.077083228682
$ws.Range("S2").Value = 0.0001074233166016483
$ws.Range("T2").Value = 0.0001074233166016484

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5203476666666668
$ws.Range("H3").Value = 1.561043
$ws.Range("I3").Value = 0.004105934376266647
$ws.Range("J3").Value = 0.004105934376266647
$ws.Range("O3").Value = 0.5375072133340244
$ws.Range("P3").Value = 0.5375072133340244
$ws.Range("Q3").Value = 4.741427711878001
$ws.Range("R3").Value = 42.67284940690201
$ws.Range("S3").Value = 0.002206969344719461
$ws.Range("T3").Value = 0.002206969344719461

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5203476666666668
$ws.Range("H4").Value = 1.561043
$ws.Range("I4").Value = 0.004105934376266647
$ws.Range("J4").Value = 0.004105934376266647
$ws.Range("M4").Value = 7.396838666666667
$ws.Range("N4").Value = 22.190516
$ws.Range("O4").Value = 0.4363298462101849
$ws.Range("P4").Value = 0.4363298462101849
$ws.Range("Q4").Value = 3.848927740909779
$ws.Range("R4").Value = 34.64034966818801
$ws.Range("S4").Value = 0.001791541714945537
$ws.Range("T4").Value = 0.001791541714945538

# Row 5
$ws.Range("I5").Value = 0.8361295370252257
$ws.Range("J5").Value = 0.8361295370252259
$ws.Range("M5").Value = 0.4435246666666666
$ws.Range("N5").Value = 1.330574
$ws.Range("O5").Value = 0.02616294045579069
$ws.Range("P5").Value = 0.02616294045579069
$ws.Range("Q5").Value = 46.99730463841622
$ws.Range("R5").Value = 422.9757417457459
$ws.Range("S5").Value = 0.02187560729051882
$ws.Range("T5").Value = 0.02187560729051882

# Row 6
$ws.Range("I6").Value = 0.8361295370252257
$ws.Range("J6").Value = 0.8361295370252259
$ws.Range("O6").Value = 0.5375072133340244
$ws.Range("P6").Value = 0.5375072133340244
$ws.Range("S6").Value = 0.4494256574326971
$ws.Range("T6").Value = 0.4494256574326971

# Row 7
$ws.Range("I7").Value = 0.8361295370252257
$ws.Range("J7").Value = 0.8361295370252259
$ws.Range("M7").Value = 7.396838666666667
$ws.Range("N7").Value = 22.190516
$ws.Range("O7").Value = 0.4363298462101849
$ws.Range("P7").Value = 0.4363298462101849
$ws.Range("Q7").Value = 783.7928897871517
$ws.Range("R7").Value = 7054.136008084365
$ws.Range("S7").Value = 0.3648282723020098
$ws.Range("T7").Value = 0.3648282723020099

# Row 8
$ws.Range("G8").Value = 20.24706
$ws.Range("H8").Value = 60.74118
$ws.Range("I8").Value = 0.1597645285985076
$ws.Range("J8").Value = 0.1597645285985076
$ws.Range("M8").Value = 0.4435246666666666
$ws.Range("N8").Value = 1.330574
$ws.Range("O8").Value = 0.02616294045579069
$ws.Range("P8").Value = 0.02616294045579069
$ws.Range("Q8").Value = 8.98007053748
$ws.Range("R8").Value = 80.82063483732
$ws.Range("S8").Value = 0.004179909848670223
$ws.Range("T8").Value = 0.004179909848670223

# Row 9
$ws.Range("G9").Value = 20.24706
$ws.Range("H9").Value = 60.74118
$ws.Range("I9").Value = 0.1597645285985076
$ws.Range("J9").Value = 0.1597645285985076
$ws.Range("O9").Value = 0.5375072133340244
$ws.Range("P9").Value = 0.5375072133340244
$ws.Range("Q9").Value = 184.49198010828
$ws.Range("R9").Value = 1660.42782097452
$ws.Range("S9").Value = 0.08587458655660786
$ws.Range("T9").Value = 0.08587458655660787

# Row 10
$ws.Range("G10").Value = 20.24706
$ws.Range("H10").Value = 60.74118
$ws.Range("I10").Value = 0.1597645285985076
$ws.Range("J10").Value = 0.1597645285985076
$ws.Range("M10").Value = 7.396838666666667
$ws.Range("N10").Value = 22.190516
$ws.Range("O10").Value = 0.4363298462101849
$ws.Range("P10").Value = 0.4363298462101849
$ws.Range("Q10").Value = 149.76423629432
$ws.Range("R10").Value = 1347.87812664888
$ws.Range("S10").Value = 0.06971003219322951
$ws.Range("T10").Value = 0.06971003219322952

Write-Host "Updated NATMI TPM values for rows 2-10"